$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$v1 = @'
test/0
'@
$ws.Range("A2").Value = $v1

$v2 = @'

def incr_list(l: list):
    """Return list with elements incremented by 1.
    >>> incr_list([1, 2, 3])
    [2, 3, 4]
    >>> incr_list([5, 3, 5, 2, 3, 3, 9, 0, 123])
    [6, 4, 6, 3, 4, 4, 10, 1, 124]
    """

'@
$ws.Range("B2").Value = $v2

$v3 = @'
['\n\ndef incr_list(l: list):\n    """Return list with elements incremented by 1.\n    >>> incr_list([1, 2, 3])\n    [2, 3, 4]\n    >>> incr_list([5, 3, 5, 2, 3, 3, 9, 0, 123])\n    [6, 4, 6, 3, 4, 4, 10, 1, 124]\n    """\n_list = [] for i in range(len(l)): if l[i] == \'y\': # y is a bit of a bit of a bit of a # bit of a character. # y is a bit of a character. # y is a bit of a character. # y is a bit of a character. # y is a bit of a character. # y is a bit of a character. # y is a bit of a character']
'@
$ws.Range("C2").Value = $v3

$v4 = @'
[True]
'@
$ws.Range("D2").Value = $v4

$v5 = @'
test/1
'@
$ws.Range("A3").Value = $v5

$v6 = @'

def is_prime(n):
    """Return true if a given number is prime, and false otherwise.
    >>> is_prime(6)
    False
    >>> is_prime(101)
    True
    >>> is_prime(11)
    True
    >>> is_prime(13441)
    True
    >>> is_prime(61)
    True
    >>> is_prime(4)
    False
    >>> is_prime(1)
    False
    """

'@
$ws.Range("B3").Value = $v6

$v7 = @'
['\n\ndef is_prime(n):\n    """Return true if a given number is prime, and false otherwise.\n    >>> is_prime(6)\n    False\n    >>> is_prime(101)\n    True\n    >>> is_prime(11)\n    True\n    >>> is_prime(13441)\n    True\n    >>> is_prime(61)\n    True\n    >>> is_prime(4)\n    False\n    >>> is_prime(1)\n    False\n    """\n = True # True means \'xxx\' == \'xxx\' == \'xxx\' == \'xxx\' == \'xxx\' == \'xxx\' == \'xxx\' == \'xxx\' == \'xxx\' == \'xxx\' == \'xxx\' == \'xxx\' == \'']
'@
$ws.Range("C3").Value = $v7

$v8 = @'
[True]
'@
$ws.Range("D3").Value = $v8

$v9 = @'
test/2
'@
$ws.Range("A4").Value = $v9

$v10 = @'

def add_elements(arr, k):
    """
    Given a non-empty array of integers arr and an integer k, return
    the sum of the elements with at most two digits from the first k elements of arr.
    Example:
        Input: arr = [111,21,3,4000,5,6,7,8,9], k = 4
        Output: 24 # sum of 21 + 3
    Constraints:
        1. 1 <= len(arr) <= 100
        2. 1 <= k <= len(arr)
    """

'@
$ws.Range("B4").Value = $v10

$v11 = @'
['\ndef add_elements(arr, k):\n    """\n    Given a non-empty array of integers arr and an integer k, return\n    the sum of the elements with at most two digits from the first k elements of arr.\n\n    Example:\n\n        Input: arr = [111,21,3,4000,5,6,7,8,9], k = 4\n        Output: 24 # sum of 21 + 3\n\n    Constraints:\n        1. 1 <= len(arr) <= 100\n        2. 1 <= k <= len(arr)\n    """\n 负ary orhood orhood orhood orhood orhood orhood orhood orhood orhood orhood or']
'@
$ws.Range("C4").Value = $v11

$v12 = @'
[False]
'@
$ws.Range("D4").Value = $v12

$v13 = @'
test/3
'@
$ws.Range("A5").Value = $v13

$v14 = @'

def solution(lst):
    """Given a non-empty list of integers, return the sum of all of the odd elements that are in even positions.
    Examples
    solution([5, 8, 7, 1]) ==> 12
    solution([3, 3, 3, 3, 3]) ==> 9
    solution([30, 13, 24, 321]) ==>0
    """

'@
$ws.Range("B5").Value = $v14

$v15 = @'
['\ndef solution(lst):\n    """Given a non-empty list of integers, return the sum of all of the odd elements that are in even positions.\n    \n\n    Examples\n    solution([5, 8, 7, 1]) ==> 12\n    solution([3, 3, 3, 3, 3]) ==> 9\n    solution([30, 13, 24, 321]) ==>0\n    """\n = 0 for index, piece in enumerate(lst): if piece == cardinal: cardinal = index return cardinal ']
'@
$ws.Range("C5").Value = $v15

$v16 = @'
[True]
'@
$ws.Range("D5").Value = $v16

$v17 = @'
test/4
'@
$ws.Range("A6").Value = $v17

$v18 = @'

def digits(n):
    """Given a positive integer n, return the product of the odd digits.
    Return 0 if all digits are even.
    For example:
    digits(1)  == 1
    digits(4)  == 0
    digits(235) == 15
    """

'@
$ws.Range("B6").Value = $v18

$v19 = @'
['\ndef digits(n):\n    """Given a positive integer n, return the product of the odd digits.\n    Return 0 if all digits are even.\n    For example:\n    digits(1)  == 1\n    digits(4)  == 0\n    digits(235) == 15\n    """\n 度 """ if n % 2 == 0: return 0 if n % 2 == 1: return 1 if n % 2 == 2: return 2 if n % 2 == 0: return 2 if n % 2 == 0: return 3 if n % 2 == 0: return 4 if n % 2 == 0: return 4 if n % 2 == 0: return 4 if n % 2 == 0: return 4 ']
'@
$ws.Range("C6").Value = $v19

$v20 = @'
[False]
'@
$ws.Range("D6").Value = $v20

$ws.Range("A7:D11").EntireRow.Delete()

Write-Host "edit complete"